$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap row 2 and row 3 (the two matches originally mis-ordered on 2023-07-22) ---
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
foreach ($col in $swapCols) {
    $rng2 = $ws.Range("$col" + "2")
    $rng3 = $ws.Range("$col" + "3")
    $val2 = $rng2.Value2
    $val3 = $rng3.Value2
    $rng2.Value2 = $val3
    $rng3.Value2 = $val2
}

# --- Step 2: append new rows 132-140 (latest round of fixtures/results) ---
# Row 132 (id 130, match 6816446)
$ws.Range("A131").Copy()
$ws.Range("A132").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E132").PasteSpecial(-4122)
$ws.Range("A132").Value2 = 130
$ws.Range("B132").Value2 = 6816446
$ws.Range("C132").Value2 = "Slovenia Prva Liga"
$ws.Range("D132").Value2 = "Slovenia Prva Liga"
$ws.Range("E132").Value2 = 45380.45833333334
$ws.Range("F132").Value2 = "NK Rogaska"
$ws.Range("G132").Value2 = "NK Bravo"
$ws.Range("H132").Value2 = 2
$ws.Range("I132").Value2 = 0
$ws.Range("J132").Value2 = "H"
$ws.Range("K132").Value2 = 2.8
$ws.Range("L132").Value2 = 3.5
$ws.Range("M132").Value2 = 2.1
$ws.Range("N132").Value2 = 2.75
$ws.Range("O132").Value2 = 3.3
$ws.Range("P132").Value2 = 2.2
$ws.Range("Q132").Value2 = 0.25
$ws.Range("R132").Value2 = 1.8
$ws.Range("S132").Value2 = 2
$ws.Range("T132").Value2 = 2.25
$ws.Range("U132").Value2 = 1.925
$ws.Range("V132").Value2 = 1.875
$ws.Range("W132").Value2 = 1.75
$ws.Range("X132").Value2 = -1
$ws.Range("Y132").Value2 = -1
$ws.Range("Z132").Value2 = 0.8
$ws.Range("AA132").Value2 = -1
$ws.Range("AB132").Value2 = -0.5
$ws.Range("AC132").Value2 = 0.4375

# Row 133 (id 131, match 7977922)
$ws.Range("A131").Copy()
$ws.Range("A133").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E133").PasteSpecial(-4122)
$ws.Range("A133").Value2 = 131
$ws.Range("B133").Value2 = 7977922
$ws.Range("C133").Value2 = "Slovenia Prva Liga"
$ws.Range("D133").Value2 = "Slovenia Prva Liga"
$ws.Range("E133").Value2 = 45380.5625
$ws.Range("F133").Value2 = "NK Maribor"
$ws.Range("G133").Value2 = "NK Radomlje"
$ws.Range("H133").Value2 = 1
$ws.Range("I133").Value2 = 0
$ws.Range("J133").Value2 = "H"
$ws.Range("K133").Value2 = 1.285
$ws.Range("L133").Value2 = 5.5
$ws.Range("M133").Value2 = 6.5
$ws.Range("N133").Value2 = 1.333
$ws.Range("O133").Value2 = 5.25
$ws.Range("P133").Value2 = 5.5
$ws.Range("Q133").Value2 = -1.5
$ws.Range("R133").Value2 = 2
$ws.Range("S133").Value2 = 1.8
$ws.Range("T133").Value2 = 2.75
$ws.Range("U133").Value2 = 1.775
$ws.Range("V133").Value2 = 2.025
$ws.Range("W133").Value2 = 0.333
$ws.Range("X133").Value2 = -1
$ws.Range("Y133").Value2 = -1
$ws.Range("Z133").Value2 = -1
$ws.Range("AA133").Value2 = 0.8
$ws.Range("AB133").Value2 = -1
$ws.Range("AC133").Value2 = 1.025

# Row 134 (id 132, match 7977924)
$ws.Range("A131").Copy()
$ws.Range("A134").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E134").PasteSpecial(-4122)
$ws.Range("A134").Value2 = 132
$ws.Range("B134").Value2 = 7977924
$ws.Range("C134").Value2 = "Slovenia Prva Liga"
$ws.Range("D134").Value2 = "Slovenia Prva Liga"
$ws.Range("E134").Value2 = 45380.67708333334
$ws.Range("F134").Value2 = "NK Domzale"
$ws.Range("G134").Value2 = "NS Mura"
$ws.Range("H134").Value2 = 3
$ws.Range("I134").Value2 = 5
$ws.Range("J134").Value2 = "A"
$ws.Range("K134").Value2 = 2
$ws.Range("L134").Value2 = 3.4
$ws.Range("M134").Value2 = 3.1
$ws.Range("N134").Value2 = 2.3
$ws.Range("O134").Value2 = 3.3
$ws.Range("P134").Value2 = 2.625
$ws.Range("Q134").Value2 = 0
$ws.Range("R134").Value2 = 1.775
$ws.Range("S134").Value2 = 2.025
$ws.Range("T134").Value2 = 2.25
$ws.Range("U134").Value2 = 1.8
$ws.Range("V134").Value2 = 2
$ws.Range("W134").Value2 = -1
$ws.Range("X134").Value2 = -1
$ws.Range("Y134").Value2 = 1.625
$ws.Range("Z134").Value2 = -1
$ws.Range("AA134").Value2 = 1.025
$ws.Range("AB134").Value2 = 0.8
$ws.Range("AC134").Value2 = -1

# Row 135 (id 133, match 7977921)
$ws.Range("A131").Copy()
$ws.Range("A135").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E135").PasteSpecial(-4122)
$ws.Range("A135").Value2 = 133
$ws.Range("B135").Value2 = 7977921
$ws.Range("C135").Value2 = "Slovenia Prva Liga"
$ws.Range("D135").Value2 = "Slovenia Prva Liga"
$ws.Range("E135").Value2 = 45381.45833333334
$ws.Range("F135").Value2 = "Olimpija Ljubljana"
$ws.Range("G135").Value2 = "FC Koper"
$ws.Range("H135").Value2 = 3
$ws.Range("I135").Value2 = 2
$ws.Range("J135").Value2 = "H"
$ws.Range("K135").Value2 = 1.5
$ws.Range("L135").Value2 = 3.6
$ws.Range("M135").Value2 = 6
$ws.Range("N135").Value2 = 1.55
$ws.Range("O135").Value2 = 3.6
$ws.Range("P135").Value2 = 5.5
$ws.Range("Q135").Value2 = -1
$ws.Range("R135").Value2 = 1.95
$ws.Range("S135").Value2 = 1.85
$ws.Range("T135").Value2 = 2.75
$ws.Range("U135").Value2 = 2
$ws.Range("V135").Value2 = 1.8
$ws.Range("W135").Value2 = 0.55
$ws.Range("X135").Value2 = -1
$ws.Range("Y135").Value2 = -1
$ws.Range("Z135").Value2 = 0
$ws.Range("AA135").Value2 = -0.0
$ws.Range("AB135").Value2 = 1
$ws.Range("AC135").Value2 = -1

# Row 136 (id 134, match 7977923)
$ws.Range("A131").Copy()
$ws.Range("A136").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E136").PasteSpecial(-4122)
$ws.Range("A136").Value2 = 134
$ws.Range("B136").Value2 = 7977923
$ws.Range("C136").Value2 = "Slovenia Prva Liga"
$ws.Range("D136").Value2 = "Slovenia Prva Liga"
$ws.Range("E136").Value2 = 45381.5625
$ws.Range("F136").Value2 = "NK Celje"
$ws.Range("G136").Value2 = "NK Aluminij"
$ws.Range("H136").Value2 = 2
$ws.Range("I136").Value2 = 2
$ws.Range("J136").Value2 = "D"
$ws.Range("K136").Value2 = 1.2
$ws.Range("L136").Value2 = 6.5
$ws.Range("M136").Value2 = 8
$ws.Range("N136").Value2 = 1.222
$ws.Range("O136").Value2 = 6.5
$ws.Range("P136").Value2 = 10
$ws.Range("Q136").Value2 = -1.75
$ws.Range("R136").Value2 = 1.8
$ws.Range("S136").Value2 = 2
$ws.Range("T136").Value2 = 3
$ws.Range("U136").Value2 = 1.85
$ws.Range("V136").Value2 = 1.95
$ws.Range("W136").Value2 = -1
$ws.Range("X136").Value2 = 5.5
$ws.Range("Y136").Value2 = -1
$ws.Range("Z136").Value2 = -1
$ws.Range("AA136").Value2 = 1
$ws.Range("AB136").Value2 = 0.8500000000000001
$ws.Range("AC136").Value2 = -1

# Row 137 (id 135, match 6814753)
$ws.Range("A131").Copy()
$ws.Range("A137").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E137").PasteSpecial(-4122)
$ws.Range("A137").Value2 = 135
$ws.Range("B137").Value2 = 6814753
$ws.Range("C137").Value2 = "Slovenia Prva Liga"
$ws.Range("D137").Value2 = "Slovenia Prva Liga"
$ws.Range("E137").Value2 = 45387.63541666666
$ws.Range("F137").Value2 = "NK Aluminij"
$ws.Range("G137").Value2 = "NK Maribor"
$ws.Range("K137").Value2 = 6
$ws.Range("L137").Value2 = 4.5
$ws.Range("M137").Value2 = 1.444
$ws.Range("N137").Value2 = 6
$ws.Range("O137").Value2 = 4.5
$ws.Range("P137").Value2 = 1.444
$ws.Range("Q137").Value2 = 1.25
$ws.Range("R137").Value2 = 1.85
$ws.Range("S137").Value2 = 1.95
$ws.Range("T137").Value2 = 2.75
$ws.Range("U137").Value2 = 1.825
$ws.Range("V137").Value2 = 1.975
$ws.Range("W137").Value2 = 0
$ws.Range("X137").Value2 = 0
$ws.Range("Y137").Value2 = 0
$ws.Range("Z137").Value2 = 0
$ws.Range("AA137").Value2 = 0

# Row 138 (id 136, match 8035687)
$ws.Range("A131").Copy()
$ws.Range("A138").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E138").PasteSpecial(-4122)
$ws.Range("A138").Value2 = 136
$ws.Range("B138").Value2 = 8035687
$ws.Range("C138").Value2 = "Slovenia Prva Liga"
$ws.Range("D138").Value2 = "Slovenia Prva Liga"
$ws.Range("E138").Value2 = 45388.41666666666
$ws.Range("F138").Value2 = "NK Rogaska"
$ws.Range("G138").Value2 = "Olimpija Ljubljana"
$ws.Range("K138").Value2 = 6
$ws.Range("L138").Value2 = 4.333
$ws.Range("M138").Value2 = 1.45
$ws.Range("N138").Value2 = 6
$ws.Range("O138").Value2 = 4.333
$ws.Range("P138").Value2 = 1.45
$ws.Range("Q138").Value2 = 1.25
$ws.Range("R138").Value2 = 1.775
$ws.Range("S138").Value2 = 2.025
$ws.Range("T138").Value2 = 3
$ws.Range("U138").Value2 = 2.025
$ws.Range("V138").Value2 = 1.775
$ws.Range("W138").Value2 = 0
$ws.Range("X138").Value2 = 0
$ws.Range("Y138").Value2 = 0
$ws.Range("Z138").Value2 = 0
$ws.Range("AA138").Value2 = 0

# Row 139 (id 137, match 6814435)
$ws.Range("A131").Copy()
$ws.Range("A139").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E139").PasteSpecial(-4122)
$ws.Range("A139").Value2 = 137
$ws.Range("B139").Value2 = 6814435
$ws.Range("C139").Value2 = "Slovenia Prva Liga"
$ws.Range("D139").Value2 = "Slovenia Prva Liga"
$ws.Range("E139").Value2 = 45388.52083333334
$ws.Range("F139").Value2 = "NK Radomlje"
$ws.Range("G139").Value2 = "FC Koper"
$ws.Range("K139").Value2 = 2.55
$ws.Range("L139").Value2 = 3.25
$ws.Range("M139").Value2 = 2.55
$ws.Range("N139").Value2 = 2.55
$ws.Range("O139").Value2 = 3.25
$ws.Range("P139").Value2 = 2.55
$ws.Range("Q139").Value2 = 0
$ws.Range("R139").Value2 = 1.9
$ws.Range("S139").Value2 = 1.9
$ws.Range("T139").Value2 = 2.25
$ws.Range("U139").Value2 = 1.775
$ws.Range("V139").Value2 = 2.025
$ws.Range("W139").Value2 = 0
$ws.Range("X139").Value2 = 0
$ws.Range("Y139").Value2 = 0
$ws.Range("Z139").Value2 = 0
$ws.Range("AA139").Value2 = 0

# Row 140 (id 138, match 6837117)
$ws.Range("A131").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("E131").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$ws.Range("A140").Value2 = 138
$ws.Range("B140").Value2 = 6837117
$ws.Range("C140").Value2 = "Slovenia Prva Liga"
$ws.Range("D140").Value2 = "Slovenia Prva Liga"
$ws.Range("E140").Value2 = 45388.63541666666
$ws.Range("F140").Value2 = "NS Mura"
$ws.Range("G140").Value2 = "NK Celje"
$ws.Range("K140").Value2 = 5.25
$ws.Range("L140").Value2 = 4.2
$ws.Range("M140").Value2 = 1.5
$ws.Range("N140").Value2 = 5.25
$ws.Range("O140").Value2 = 4.2
$ws.Range("P140").Value2 = 1.5
$ws.Range("Q140").Value2 = 1
$ws.Range("R140").Value2 = 1.95
$ws.Range("S140").Value2 = 1.85
$ws.Range("T140").Value2 = 2.75
$ws.Range("U140").Value2 = 1.975
$ws.Range("V140").Value2 = 1.825
$ws.Range("W140").Value2 = 0
$ws.Range("X140").Value2 = 0
$ws.Range("Y140").Value2 = 0
$ws.Range("Z140").Value2 = 0
$ws.Range("AA140").Value2 = 0

$excel.CutCopyMode = $false